$d = $word.ActiveDocument

# Locate the title paragraph ("Dheeraj Chand") and work off its own
# paragraph Range, so the insertion point is anchored to content rather
# than a hard-coded paragraph index.
$titleRange = $d.Content
$titleRange.Find.Execute("Dheeraj Chand", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$titlePara = $titleRange.Paragraphs.Item(1)

# Insert a brand-new (empty) paragraph right after the title. It comes
# back centered (inherited from the title's paragraph formatting), which
# is exactly the alignment the contact line needs.
$titlePara.Range.InsertParagraphAfter()
$contactPara = $d.Paragraphs.Item($titlePara.Index + 1)

# Populate the new paragraph purely via OOXML so the run gets no
# rPr/formatting baggage inherited from the bold/28pt title run - just a
# plain run of text, centered by its own pPr.
$contactXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r>' + `
  '<w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX</w:t>' + `
  '</w:r></w:p></w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$contactPara.Range.InsertXML($contactXml)
